$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.298.93'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '2.526.56'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Value = '2.525.58'
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.34%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.349'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '2.959.46'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.32'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '59.188.30'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '2.525.09'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.426'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0781'
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.21%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  +2.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.75'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.828'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '285.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.16%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.608'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0935'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0512'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.85%  '
